$d = $word.ActiveDocument

# --- Step 1: grow paragraph count to 27, reusing the first 17 existing paragraphs ---
$targetCount = 27
while ($d.Paragraphs.Count -lt $targetCount) {
    $last = $d.Paragraphs($d.Paragraphs.Count)
    $last.Range.InsertParagraphAfter()
}

# --- Step 2: set the text of every paragraph to the new (reordered / trimmed) content ---
$d.Paragraphs(1).Range.Text = "command-line application that accepts user input"
$d.Paragraphs(2).Range.Text = "prompted for my team members and their information"
$d.Paragraphs(3).Range.Text = "HTML file is generated that displays a nicely formatted team roster based on user input"
$d.Paragraphs(4).Range.Text = "start the application"
$d.Paragraphs(5).Range.Text = "I am prompted to enter the team manager’s name, employee ID, email address, and office number"
$d.Paragraphs(6).Range.Text = "enter the team manager’s name, employee ID, email address, and office number"
$d.Paragraphs(7).Range.Text = "presented with a menu with the option to add an engineer or an intern or to finish building my team"
$d.Paragraphs(8).Range.Text = "I select the engineer option"
$d.Paragraphs(9).Range.Text = "prompted to enter the engineer’s name, ID, email, and GitHub username, and I am taken back to the menu"
$d.Paragraphs(10).Range.Text = "I select the intern option"
$d.Paragraphs(11).Range.Text = "prompted to enter the intern’s name, ID, email, and school, and I am taken back to the menu"
$d.Paragraphs(12).Range.Text = "to finish building my team"
$d.Paragraphs(13).Range.Text = "exit the application, and the HTML is generated"
$d.Paragraphs(14).Range.Text = "click on an email address in the HTML"
$d.Paragraphs(15).Range.Text = "my default email program opens and populates the TO field of the email with the address"
$d.Paragraphs(16).Range.Text = "click on the GitHub username"
$d.Paragraphs(17).Range.Text = "GitHub profile opens in a new tab"
$d.Paragraphs(18).Range.Text = "start the application"
$d.Paragraphs(19).Range.Text = "I am prompted to enter the team manager’s name, employee ID, email address, and office number"
$d.Paragraphs(20).Range.Text = "enter the team manager’s name, employee ID, email address, and office number"
$d.Paragraphs(21).Range.Text = "presented with a menu with the option to add an engineer or an intern or to finish building my team"
$d.Paragraphs(22).Range.Text = "I select the engineer option"
$d.Paragraphs(23).Range.Text = "prompted to enter the engineer’s name, ID, email, and GitHub username, and I am taken back to the menu"
$d.Paragraphs(24).Range.Text = "I select the intern option"
$d.Paragraphs(25).Range.Text = "prompted to enter the intern’s name, ID, email, and school, and I am taken back to the menu"
$d.Paragraphs(26).Range.Text = "to finish building my team"
$d.Paragraphs(27).Range.Text = "exit the application, and the HTML is generated"

# --- Step 3: turn the whole block into one shared bulleted list (List Paragraph style, numId 1) ---
$listRange = $d.Range($d.Paragraphs(1).Range.Start, $d.Paragraphs($d.Paragraphs.Count).Range.End)
$listRange.Style = "List Paragraph"
$listRange.ListFormat.ApplyBulletDefault()

# --- Step 4: set the correct outline level (ilvl) on every paragraph ---
$d.Paragraphs(2).Range.ListFormat.ListLevelNumber = 2
$d.Paragraphs(3).Range.ListFormat.ListLevelNumber = 2
$d.Paragraphs(5).Range.ListFormat.ListLevelNumber = 2
$d.Paragraphs(6).Range.ListFormat.ListLevelNumber = 2
$d.Paragraphs(7).Range.ListFormat.ListLevelNumber = 3
$d.Paragraphs(8).Range.ListFormat.ListLevelNumber = 3
$d.Paragraphs(9).Range.ListFormat.ListLevelNumber = 4
$d.Paragraphs(10).Range.ListFormat.ListLevelNumber = 3
$d.Paragraphs(11).Range.ListFormat.ListLevelNumber = 4
$d.Paragraphs(12).Range.ListFormat.ListLevelNumber = 3
$d.Paragraphs(13).Range.ListFormat.ListLevelNumber = 4
$d.Paragraphs(15).Range.ListFormat.ListLevelNumber = 2
$d.Paragraphs(17).Range.ListFormat.ListLevelNumber = 2
$d.Paragraphs(19).Range.ListFormat.ListLevelNumber = 2
$d.Paragraphs(20).Range.ListFormat.ListLevelNumber = 2
$d.Paragraphs(21).Range.ListFormat.ListLevelNumber = 3
$d.Paragraphs(22).Range.ListFormat.ListLevelNumber = 3
$d.Paragraphs(23).Range.ListFormat.ListLevelNumber = 4
$d.Paragraphs(24).Range.ListFormat.ListLevelNumber = 3
$d.Paragraphs(25).Range.ListFormat.ListLevelNumber = 4
$d.Paragraphs(26).Range.ListFormat.ListLevelNumber = 3
$d.Paragraphs(27).Range.ListFormat.ListLevelNumber = 4

# --- Step 5: match the "List Paragraph" style definition (uiPriority/indent/contextualSpacing) ---
$s = $d.Styles("List Paragraph")
$s.Priority = 34
$s.ParagraphFormat.LeftIndent = 36
$s.NoSpaceBetweenParagraphsOfSameStyle = $true

# --- Step 6: split "I select the intern option" into two runs ("...optio" + "n") ---
$splitPara = $d.Paragraphs(24)
$splitRange = $splitPara.Range
$endPos = $splitRange.End
$lastCharRange = $d.Range($endPos - 2, $endPos - 1)
$lastCharRange.Bold = 1
$lastCharRange.Bold = 0
